$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "import" row paired with the existing A1 import row (row 13)
$ws.Range("B14").Value = "import"
$ws.Range("C14").Value = "org.openl.rules.beans.B1"

# New method rows describing the two B1 "hello" overloads
$ws.Range("B21").Value = "Method B1 hello1()"
$ws.Range("B22").Value = "return B1(name=""hello"");"

$ws.Range("B26").Value = "Method B1 hello2()"
$ws.Range("B27").Value = "return B1(var=""hello"");"

# Move the selection like the authored workbook
$ws.Range("C9").Select()
